$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "0.9982", "47.50", "24.913.19") are preserved exactly as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.017.17"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "1.664.15"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.78%  "

$ws.Range("D5").Value = "325.61"
$ws.Range("E5").Value = "  +5.29%  "

$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("D7").Value = "0.3641"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").Value = "47.50"
$ws.Range("E8").Value = "  +0.80%  "

$ws.Range("D9").Value = "0.3277"
$ws.Range("E9").Value = "  -1.48%  "

$ws.Range("D10").Value = "1.140"
$ws.Range("E10").Value = "  -1.10%  "

$ws.Range("D11").Value = "0.07101"
$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("D12").Value = "0.9956"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").Value = "6.072"
$ws.Range("E13").Value = "  -0.71%  "

$ws.Range("D14").Value = "19.68"
$ws.Range("E14").Value = "  -2.42%  "

$ws.Range("D15").Value = "1.662.85"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").Value = "6.623"
$ws.Range("E16").Value = "  -1.69%  "

$ws.Range("D17").Value = "0.00001052"
$ws.Range("E17").Value = "  -2.75%  "

$ws.Range("D18").Value = "0.06604"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").Value = "0.9969"
$ws.Range("E19").Value = "  -0.81%  "

$ws.Range("D20").Value = "79.38"
$ws.Range("E20").Value = "  -2.01%  "

$ws.Range("D21").Value = "5.933"
$ws.Range("E21").Value = "  -2.18%  "

$ws.Range("D22").Value = "15.85"
$ws.Range("E22").Value = "  -4.33%  "

$ws.Range("D23").Value = "12.70"
$ws.Range("E23").Value = "  +2.55%  "

$ws.Range("D24").Value = "24.933.22"
$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("D25").Value = "2.440"
$ws.Range("E25").Value = "  +1.15%  "

$ws.Range("D26").Value = "2.460"
$ws.Range("E26").Value = "  -6.10%  "

$ws.Range("D27").Value = "147.96"
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").Value = "18.67"
$ws.Range("E28").Value = "  -4.62%  "

$ws.Range("D29").Value = "1.841.31"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "1.207"
$ws.Range("E30").Value = "  +0.59%  "

$ws.Range("D31").Value = "125.23"
$ws.Range("E31").Value = "  -2.88%  "

$ws.Range("D32").Value = "4.095"
$ws.Range("E32").Value = "  -1.04%  "

$ws.Range("D33").Value = "5.802"
$ws.Range("E33").Value = "  -7.81%  "

$ws.Range("D34").Value = "0.08459"
$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("E35").Value = "  -5.37%  "

$ws.Range("D36").Value = "12.32"
$ws.Range("E36").Value = "  -5.85%  "

$ws.Range("D37").Value = "1.295"
$ws.Range("E37").Value = "  +5.18%  "

$ws.Range("D38").Value = "5.200"
$ws.Range("E38").Value = "  -2.18%  "

$ws.Range("D39").Value = "0.02282"
$ws.Range("E39").Value = "  -1.52%  "

$ws.Range("D40").Value = "0.06099"
$ws.Range("E40").Value = "  -3.42%  "

$ws.Range("D41").Value = "8.386"
$ws.Range("E41").Value = "  -2.85%  "

$ws.Range("D42").Value = "0.2079"
$ws.Range("E42").Value = "  -2.19%  "

$ws.Range("D43").Value = "0.9968"
$ws.Range("E43").Value = "  -0.58%  "

$ws.Range("D44").Value = "0.5964"
$ws.Range("E44").Value = "  -3.36%  "

$ws.Range("D45").Value = "13.91"
$ws.Range("E45").Value = "  +5.25%  "

$ws.Range("D46").Value = "3.859"
$ws.Range("E46").Value = "  +2.07%  "

$ws.Range("D47").Value = "0.5652"
$ws.Range("E47").Value = "  -3.69%  "

$ws.Range("D48").Value = "125.40"
$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").Value = "1.959"
$ws.Range("E49").Value = "  -2.70%  "

$ws.Range("D50").Value = "0.07001"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("D51").Value = "1.192"
$ws.Range("E51").Value = "  -0.30%  "
